$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (prices, 1h volume %, and
# two row swaps: WrappedEther/TRON and Litecoin/Chainlink and EnergySwap/Quant)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.848.63'
$ws.Range('E2').Value = '  -2.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.815.93'
$ws.Range('E4').Value = '  -0.55%  '
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4617'
$ws.Range('E7').Value = '  -2.84%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3643'
$ws.Range('E8').Value = '  -1.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07216'
$ws.Range('E9').Value = '  -3.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8555'
$ws.Range('E10').Value = '  -3.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.69'
$ws.Range('E11').Value = '  -3.83%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07538'
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.839.46'
$ws.Range('E13').Value = '  -2.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.319'
$ws.Range('E14').Value = '  -2.52%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.493'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.61'
$ws.Range('E16').Value = '  -1.80%  '
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('E18').Value = '  -2.67%  '
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.022.13'
$ws.Range('E20').Value = '  -1.79%  '
$ws.Range('E21').Value = '  -2.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.142'
$ws.Range('E22').Value = '  -3.55%  '
$ws.Range('E23').Value = '  -2.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.071.72'
$ws.Range('E24').Value = '  -1.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.43'
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.837'
$ws.Range('E26').Value = '  -3.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.05'
$ws.Range('E27').Value = '  -3.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.056'
$ws.Range('E28').Value = '  -4.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.076'
$ws.Range('E29').Value = '  -3.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '114.61'
$ws.Range('E30').Value = '  -2.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08854'
$ws.Range('E31').Value = '  -1.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.958'
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('E33').Value = '  -4.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.390'
$ws.Range('E34').Value = '  -3.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7182'
$ws.Range('E35').Value = '  -5.28%  '
$ws.Range('E36').Value = '  -0.59%  '
$ws.Range('E37').Value = '  -3.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05231'
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.420'
$ws.Range('E39').Value = '  +0.70%  '
$ws.Range('E40').Value = '  -2.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.915'
$ws.Range('E41').Value = '  -2.96%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.122'
$ws.Range('E42').Value = '  -2.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5128'
$ws.Range('E43').Value = '  -4.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1618'
$ws.Range('E44').Value = '  -2.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.153'
$ws.Range('E45').Value = '  -4.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4783'
$ws.Range('E46').Value = '  -2.90%  '
$ws.Range('E47').Value = '  -0.58%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.06'
$ws.Range('E48').Value = '  -4.59%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.86'
$ws.Range('E49').Value = '  -1.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.611'
$ws.Range('E50').Value = '  -4.21%  '
$ws.Range('E51').Value = '  -1.97%  '
